$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hoang Viet Bach"

$ws.Range("G12").Select()
